# Update the "carjacking by month / year-over-year" workbook with the
# data refresh for 2021-12-17 (report now covers "through 12-09" instead
# of "through 12-08").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet title reflects the new as-of date.
$ws.Name = "Through 2021-12-09"

# Row label for December also reflects the new as-of date.
$ws.Range("A13").Value = "December (through 12-09)"

# 2021 column, November row: one more incident.
$ws.Range("H12").Value = 201

# December row (partial month) picked up additional incidents across
# several years.
$ws.Range("C13").Value = 26
$ws.Range("D13").Value = 36
$ws.Range("E13").Value = 21
$ws.Range("F13").Value = 13
$ws.Range("G13").Value = 42
$ws.Range("H13").Value = 71

# Totals row updated to match the new December figures.
$ws.Range("C14").Value = 589
$ws.Range("D14").Value = 857
$ws.Range("E14").Value = 703
$ws.Range("F14").Value = 547
$ws.Range("G14").Value = 1306
$ws.Range("H14").Value = 1714
